$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D text values (e.g. "61.943.13", "0.0520") are not
# auto-coerced to numbers by Excel when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.943.13"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "2.434.31"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "578.91"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").Value = "142.32"
$ws.Range("E6").Value = "  -3.10%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").Value = "2.431.52"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").Value = "0.106"
$ws.Range("E10").Value = "  -4.00%  "

$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  -3.52%  "

$ws.Range("D14").Value = "26.26"
$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.874.50"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -4.68%  "

$ws.Range("D17").Value = "62.078.44"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "2.428.36"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").Value = "10.86"
$ws.Range("E19").Value = "  -4.11%  "

$ws.Range("D20").Value = "7.04"
$ws.Range("E20").Value = "  -4.44%  "

$ws.Range("D21").Value = "328.69"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -2.45%  "

$ws.Range("D23").Value = "1.94"
$ws.Range("E23").Value = "  -6.51%  "

$ws.Range("E24").Value = "  -2.64%  "

$ws.Range("D25").Value = "65.56"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").Value = "9.26"
$ws.Range("E26").Value = "  +4.38%  "

$ws.Range("D27").Value = "611.91"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "2.555.13"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  -9.03%  "

$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -6.78%  "

$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -4.21%  "

$ws.Range("D33").Value = "0.141"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("D35").Value = "4.87"
$ws.Range("E35").Value = "  -6.44%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("E37").Value = "  -7.19%  "

$ws.Range("D38").Value = "0.373"
$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("D39").Value = "151.09"
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("D40").Value = "18.26"
$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("D41").Value = "5.20"
$ws.Range("E41").Value = "  -4.45%  "

$ws.Range("D42").Value = "1.74"
$ws.Range("E42").Value = "  -3.14%  "

$ws.Range("D43").Value = "42.76"
$ws.Range("E43").Value = "  +2.13%  "

$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  -10.08%  "

$ws.Range("D46").Value = "142.21"
$ws.Range("E46").Value = "  -4.55%  "

$ws.Range("D47").Value = "3.59"
$ws.Range("E47").Value = "  -4.24%  "

$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  -3.02%  "

$ws.Range("D49").Value = "0.596"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "19.38"
$ws.Range("E50").Value = "  -9.24%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.0901"
$ws.Range("E51").Value = "  -1.64%  "

# Restore default "Normal" style on column D so only the value changed
$ws.Range("D2:D51").Style = "Normal"
